$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format before writing numeric-looking values,
# then clear the number-format again so the cells fall back to the default
# style (matches the source workbook, which has no explicit style on data rows).
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Range("D2").Value = '70.960.24'
$ws.Range("E2").Value = '  -2.55%  '

$ws.Range("D3").Value = '3.867.19'
$ws.Range("E3").Value = '  -2.69%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").Value = '594.85'
$ws.Range("E5").Value = '  +1.15%  '

$ws.Range("D6").Value = '166.45'
$ws.Range("E6").Value = '  +4.72%  '

$ws.Range("D7").Value = '0.673'
$ws.Range("E7").Value = '  -1.74%  '

$ws.Range("E8").Value = '  +0.11%  '

$ws.Range("D9").Value = '0.754'
$ws.Range("E9").Value = '  +0.83%  '

$ws.Range("D10").Value = '0.176'
$ws.Range("E10").Value = '  +4.83%  '

$ws.Range("D11").Value = '53.43'
$ws.Range("E11").Value = '  -1.18%  '

$ws.Range("E12").Value = '  +1.17%  '

$ws.Range("D13").Value = '11.45'
$ws.Range("E13").Value = '  +5.38%  '

$ws.Range("D14").Value = '4.488.30'
$ws.Range("E14").Value = '  -2.70%  '

$ws.Range("D15").Value = '21.43'
$ws.Range("E15").Value = '  +5.47%  '

$ws.Range("D16").Value = '3.872.58'
$ws.Range("E16").Value = '  -2.79%  '

$ws.Range("D17").Value = '13.86'
$ws.Range("E17").Value = '  -1.23%  '

$ws.Range("E18").Value = '  -4.07%  '

$ws.Range("E19").Value = '  -2.14%  '

$ws.Range("D20").Value = '70.845.58'
$ws.Range("E20").Value = '  -2.47%  '

$ws.Range("D21").Value = '437.40'
$ws.Range("E21").Value = '  +0.71%  '

$ws.Range("D22").Value = '4.72'
$ws.Range("E22").Value = '  +0.25%  '

$ws.Range("D23").Value = '94.29'
$ws.Range("E23").Value = '  -1.75%  '

$ws.Range("D24").Value = '3.27'
$ws.Range("E24").Value = '  -4.28%  '

$ws.Range("D25").Value = '13.92'
$ws.Range("E25").Value = '  -2.59%  '

$ws.Range("D26").Value = '11.33'
$ws.Range("E26").Value = '  +1.48%  '

$ws.Range("E27").Value = '  -8.71%  '

$ws.Range("B28").Value = 'Filecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D28").Value = '10.40'
$ws.Range("E28").Value = '  -1.63%  '

$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '35.18'
$ws.Range("E29").Value = '  -3.08%  '

$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").Value = '8.07'
$ws.Range("E30").Value = '  +2.97%  '

$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").Value = '13.58'
$ws.Range("E31").Value = '  -0.45%  '

$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").Value = '48.32'
$ws.Range("E32").Value = '  -0.11%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.125'
$ws.Range("E33").Value = '  -4.18%  '

$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").Value = '69.88'
$ws.Range("E34").Value = '  +0.34%  '

$ws.Range("B35").Value = 'PEPE'
$ws.Range("C35").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D35").Value = '0.0₃0986'
$ws.Range("E35").Value = '  +12.95%  '

$ws.Range("B36").Value = 'Bittensor'
$ws.Range("C36").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D36").Value = '636.25'
$ws.Range("E36").Value = '  -6.69%  '

$ws.Range("B37").Value = 'TheGraph'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D37").Value = '0.425'
$ws.Range("E37").Value = '  -2.12%  '

$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = '0.146'
$ws.Range("E38").Value = '  -0.34%  '

$ws.Range("B39").Value = 'Dai'
$ws.Range("C39").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  +0.14%  '

$ws.Range("B40").Value = 'ThetaToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D40").Value = '3.32'
$ws.Range("E40").Value = '  -2.96%  '

$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.32%  '

$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").Value = '3.25'
$ws.Range("E42").Value = '  +25.74%  '

$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = '0.0471'
$ws.Range("E43").Value = '  -3.08%  '

$ws.Range("B44").Value = 'THORChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D44").Value = '10.09'
$ws.Range("E44").Value = '  -7.26%  '

$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").Value = '2.72'
$ws.Range("E45").Value = '  +3.20%  '

$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Value = '0.144'
$ws.Range("E46").Value = '  -3.80%  '

$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").Value = '2.85'
$ws.Range("E47").Value = '  -14.60%  '

$ws.Range("D48").Value = '3.29'
$ws.Range("E48").Value = '  -3.30%  '

$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '2.829.66'
$ws.Range("E49").Value = '  +0.90%  '

$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").Value = '0.000272'
$ws.Range("E50").Value = '  +1.78%  '

$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").Value = '2.91'
$ws.Range("E51").Value = '  -3.47%  '

$colD.ClearFormats()